$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 16: new worker DANY ALFONSO BRAVO JULIO (73007298), periodo 2306
$ws.Range("C16").Value = "73007298"
$ws.Range("D16").Value = "DANY ALFONSO BRAVO JULIO"
$ws.Range("E16").Value = "2306"
$ws.Range("F16").Value = 46400
$ws.Range("G16").Value = 1160000

# Row 17: MARIA INELDA BLANCO MORALES (64917008), periodo 2305
$ws.Range("C17").Value = "64917008"
$ws.Range("D17").Value = "MARIA INELDA BLANCO MORALES"
$ws.Range("E17").Value = "2305"
$ws.Range("F17").Value = 46400
$ws.Range("G17").Value = 1160000

# Row 18: MARIA INELDA BLANCO MORALES (64917008), periodo 2304
$ws.Range("C18").Value = "64917008"
$ws.Range("D18").Value = "MARIA INELDA BLANCO MORALES"
$ws.Range("E18").Value = "2304"
$ws.Range("F18").Value = 46400
$ws.Range("G18").Value = 1160000

# Row 19: MARIA INELDA BLANCO MORALES (64917008), periodo 2301
$ws.Range("C19").Value = "64917008"
$ws.Range("D19").Value = "MARIA INELDA BLANCO MORALES"
$ws.Range("E19").Value = "2301"
$ws.Range("F19").Value = 40000
$ws.Range("G19").Value = 1160000
